# Remove the "Nicolás" and "Catary" individual-conclusion paragraphs
# (their headings plus the English and Spanish reflection text blocks)
# from the "10. Conclusiones Individuales" section, leaving only the
# "Cristian" entries followed directly by section "11. Reflexión".

$d = $word.ActiveDocument

$startText = "Nicolás:"
$endText = "El enfoque estructurado proporcionado"

$count = $d.Paragraphs.Count
$startIdx = -1
$endIdx = -1

for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($startIdx -eq -1 -and $t.StartsWith($startText)) {
        $startIdx = $i
    }
    if ($t.StartsWith($endText)) {
        $endIdx = $i
    }
}

if ($startIdx -ne -1 -and $endIdx -ne -1) {
    $rangeStart = $d.Paragraphs.Item($startIdx).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIdx).Range.End
    $deleteRange = $d.Range($rangeStart, $rangeEnd)
    $deleteRange.Delete()
}
